$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 13956.5
$ws.Range("I40").Value = 15664.571
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 15664.571
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -15489.571
$ws.Range("N40").Value = -2350
$ws.Range("H64").Value = 7859148.5
$ws.Range("I64").Value = 9261993
$ws.Range("J64").Value = 5053460.5
$ws.Range("K64").Value = 9261993
$ws.Range("L64").Value = 5053460.5
$ws.Range("M64").Value = -9261745
$ws.Range("N64").Value = -5053956.5
$ws.Range("H67").Value = 7859148.5
$ws.Range("I67").Value = 9261993
$ws.Range("J67").Value = 5053460.5
$ws.Range("K67").Value = 9261993
$ws.Range("L67").Value = 5053460.5
$ws.Range("M67").Value = -9261135
$ws.Range("N67").Value = -5055176.5
$ws.Range("H96").Value = 20344.889
$ws.Range("I96").Value = 6188.6665
$ws.Range("J96").Value = 48657.332
$ws.Range("K96").Value = 18565.9995
$ws.Range("L96").Value = 145971.996
$ws.Range("M96").Value = -17192.9995
$ws.Range("N96").Value = -148717.996
$ws.Range("H107").Value = 4942.478
$ws.Range("I107").Value = 5184.619
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 5184.619
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -3264.619
$ws.Range("N107").Value = -6240
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 19061.666
$ws.Range("I63").Value = 2451.111
$ws.Range("J63").Value = 68893.336
$ws.Range("K63").Value = 2451.111
$ws.Range("L63").Value = 68893.336
$ws.Range("M63").Value = -1765.111
$ws.Range("N63").Value = -70265.336
$ws.Range("H66").Value = 19061.666
$ws.Range("I66").Value = 2451.111
$ws.Range("J66").Value = 68893.336
$ws.Range("K66").Value = 12255.555
$ws.Range("L66").Value = 344466.68
$ws.Range("M66").Value = -8823.555
$ws.Range("N66").Value = -351330.68
$ws.Range("H102").Value = 17316.924
$ws.Range("I102").Value = 2301.6667
$ws.Range("J102").Value = 197500
$ws.Range("K102").Value = 2301.6667
$ws.Range("L102").Value = 197500
$ws.Range("M102").Value = -679.6667000000002
$ws.Range("N102").Value = -200744
$ws.Range("H122").Value = 1546.762
$ws.Range("I122").Value = 1429.7059
$ws.Range("J122").Value = 2044.25
$ws.Range("K122").Value = 4289.1177
$ws.Range("L122").Value = 6132.75
$ws.Range("M122").Value = -1839.1177
$ws.Range("N122").Value = -11032.75
$ws.Range("H132").Value = 2245
$ws.Range("I132").Value = 1684.1777
$ws.Range("K132").Value = 5052.5331
$ws.Range("M132").Value = -2522.5331
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3032.0667
$ws.Range("I105").Value = 2666.6667
$ws.Range("J105").Value = 3580.1667
$ws.Range("K105").Value = 2666.6667
$ws.Range("L105").Value = 3580.1667
$ws.Range("M105").Value = -919.6667000000002
$ws.Range("N105").Value = -7074.1667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3636.818
$ws.Range("I62").Value = 4401.25
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 4401.25
$ws.Range("L62").Value = 3200
$ws.Range("M62").Value = -3777.25
$ws.Range("N62").Value = -4448
$ws.Range("H65").Value = 3636.818
$ws.Range("I65").Value = 4401.25
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 22006.25
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = -18886.25
$ws.Range("N65").Value = -22240
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7485.3184
$ws.Range("I80").Value = 6025.1816
$ws.Range("J80").Value = 8945.454
$ws.Range("K80").Value = 6025.1816
$ws.Range("L80").Value = 8945.454
$ws.Range("M80").Value = -5027.1816
$ws.Range("N80").Value = -10941.454
$ws.Range("H83").Value = 7485.3184
$ws.Range("I83").Value = 6025.1816
$ws.Range("J83").Value = 8945.454
$ws.Range("K83").Value = 30125.908
$ws.Range("L83").Value = 44727.27
$ws.Range("M83").Value = -25133.908
$ws.Range("N83").Value = -54711.27
$ws.Range("H102").Value = 1306.4706
$ws.Range("I102").Value = 1243.5714
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 1243.5714
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = 378.4286
$ws.Range("N102").Value = -4844
$ws.Range("H135").Value = 70819.75
$ws.Range("J135").Value = 70819.75
$ws.Range("L135").Value = 70819.75
$ws.Range("N135").Value = -80959.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5717.5835
$ws.Range("I40").Value = 5326.125
$ws.Range("J40").Value = 6500.5
$ws.Range("K40").Value = 5326.125
$ws.Range("L40").Value = 6500.5
$ws.Range("M40").Value = -5190.125
$ws.Range("N40").Value = -6772.5
$ws.Range("H46").Value = 2636.2307
$ws.Range("I46").Value = 1024.5555
$ws.Range("J46").Value = 6262.5
$ws.Range("K46").Value = 1024.5555
$ws.Range("L46").Value = 6262.5
$ws.Range("M46").Value = -836.5554999999999
$ws.Range("N46").Value = -6638.5
$ws.Range("H122").Value = 37768.035
$ws.Range("I122").Value = 42060.32
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 126180.96
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -123730.96
$ws.Range("N122").Value = -10897
$ws.Range("H132").Value = 4482.6523
$ws.Range("I132").Value = 3308.6924
$ws.Range("K132").Value = 9926.0772
$ws.Range("M132").Value = -7396.0772
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 689.5
$ws.Range("I113").Value = 502.66666
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1507.99998
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = 662.0000199999999
$ws.Range("N113").Value = -8090
$ws.Range("H122").Value = 1183.25
$ws.Range("I122").Value = 1222.1111
$ws.Range("J122").Value = 1066.6666
$ws.Range("K122").Value = 3666.3333
$ws.Range("L122").Value = 3199.9998
$ws.Range("M122").Value = -1216.3333
$ws.Range("N122").Value = -8099.9998
$ws.Range("H126").Value = 2452871.5
$ws.Range("I126").Value = 3678207.5
$ws.Range("J126").Value = 2199.75
$ws.Range("K126").Value = 11034622.5
$ws.Range("L126").Value = 6599.25
$ws.Range("M126").Value = -11032152.5
$ws.Range("N126").Value = -11539.25
$ws.Range("H132").Value = 1501272.9
$ws.Range("I132").Value = 2558889.2
$ws.Range("J132").Value = 2982.9167
$ws.Range("K132").Value = 7676667.600000001
$ws.Range("L132").Value = 8948.750100000001
$ws.Range("M132").Value = -7674137.600000001
$ws.Range("N132").Value = -14008.7501
